$d = $word.ActiveDocument

$d.Content.Find.Execute("2025-02-16 Sunday", $true, $false, $false, $false, $false, $true, 1, $false, "2025-02-17 Monday", 2) | Out-Null
$d.Content.Find.Execute("516÷3=", $true, $false, $false, $false, $false, $true, 1, $false, "954÷2=", 2) | Out-Null
$d.Content.Find.Execute("625÷8=", $true, $false, $false, $false, $false, $true, 1, $false, "976÷4=", 2) | Out-Null
$d.Content.Find.Execute("647÷8=", $true, $false, $false, $false, $false, $true, 1, $false, "819÷7=", 2) | Out-Null
$d.Content.Find.Execute("887÷9=", $true, $false, $false, $false, $false, $true, 1, $false, "552÷9=", 2) | Out-Null
$d.Content.Find.Execute("196÷7=", $true, $false, $false, $false, $false, $true, 1, $false, "826÷4=", 2) | Out-Null
$d.Content.Find.Execute("332÷6=", $true, $false, $false, $false, $false, $true, 1, $false, "481÷7=", 2) | Out-Null
$d.Content.Find.Execute("802÷6=", $true, $false, $false, $false, $false, $true, 1, $false, "469÷4=", 2) | Out-Null
$d.Content.Find.Execute("721÷5=", $true, $false, $false, $false, $false, $true, 1, $false, "374÷7=", 2) | Out-Null
$d.Content.Find.Execute("778÷3=", $true, $false, $false, $false, $false, $true, 1, $false, "590÷9=", 2) | Out-Null
$d.Content.Find.Execute("497÷5=", $true, $false, $false, $false, $false, $true, 1, $false, "429÷3=", 2) | Out-Null
$d.Content.Find.Execute("711÷6=", $true, $false, $false, $false, $false, $true, 1, $false, "816÷5=", 2) | Out-Null
$d.Content.Find.Execute("255÷6=", $true, $false, $false, $false, $false, $true, 1, $false, "241÷6=", 2) | Out-Null
$d.Content.Find.Execute("325÷4=", $true, $false, $false, $false, $false, $true, 1, $false, "210÷6=", 2) | Out-Null
$d.Content.Find.Execute("569÷5=", $true, $false, $false, $false, $false, $true, 1, $false, "517÷9=", 2) | Out-Null
$d.Content.Find.Execute("920÷2=", $true, $false, $false, $false, $false, $true, 1, $false, "573÷3=", 2) | Out-Null
$d.Content.Find.Execute("378÷2=", $true, $false, $false, $false, $false, $true, 1, $false, "966÷6=", 2) | Out-Null
$d.Content.Find.Execute("873÷4=", $true, $false, $false, $false, $false, $true, 1, $false, "435÷3=", 2) | Out-Null
$d.Content.Find.Execute("420÷2=", $true, $false, $false, $false, $false, $true, 1, $false, "290÷5=", 2) | Out-Null
$d.Content.Find.Execute("696÷3=", $true, $false, $false, $false, $false, $true, 1, $false, "651÷3=", 2) | Out-Null
$d.Content.Find.Execute("738÷3=", $true, $false, $false, $false, $false, $true, 1, $false, "422÷5=", 2) | Out-Null
$d.Content.Find.Execute("227÷5=", $true, $false, $false, $false, $false, $true, 1, $false, "433÷3=", 2) | Out-Null
$d.Content.Find.Execute("384÷3=", $true, $false, $false, $false, $false, $true, 1, $false, "984÷9=", 2) | Out-Null
$d.Content.Find.Execute("950÷8=", $true, $false, $false, $false, $false, $true, 1, $false, "545÷3=", 2) | Out-Null
$d.Content.Find.Execute("908÷4=", $true, $false, $false, $false, $false, $true, 1, $false, "969÷8=", 2) | Out-Null
$d.Content.Find.Execute("736÷9=", $true, $false, $false, $false, $false, $true, 1, $false, "979÷7=", 2) | Out-Null
